$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report (row) arrived for "Terminal Hortofrutícola Agro
# Chillán" / "Haba". The previous latest row (row 8) becomes historical data
# and is pushed down to a new row 9, while row 8 is overwritten with the
# fresh weekly figures (new date, new min/max/avg prices, new origin and
# $/Kg price).

# 1) Preserve the current (old) row 8 values by copying them down to row 9.
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = $ws.Range("B8").Value2
$ws.Range("C9").Value = $ws.Range("C8").Value2
$ws.Range("D9").Value = $ws.Range("D8").Value2
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E9").Value = $ws.Range("E8").Value2
$ws.Range("F9").Value = $ws.Range("F8").Value2
$ws.Range("G9").Value = $ws.Range("G8").Value2
$ws.Range("H9").Value = $ws.Range("H8").Value2
$ws.Range("I9").Value = $ws.Range("I8").Value2
$ws.Range("J9").Value = $ws.Range("J8").Value2
$ws.Range("K9").Value = $ws.Range("K8").Value2
$ws.Range("L9").Value = $ws.Range("L8").Value2
$ws.Range("M9").Value = $ws.Range("M8").Value2
$ws.Range("N9").Value = $ws.Range("N8").Value2
$ws.Range("O9").Value = $ws.Range("O8").Value2
$ws.Range("P9").Value = $ws.Range("P8").Value2
$ws.Range("Q9").Value = $ws.Range("Q8").Value2
$ws.Range("R9").Value = $ws.Range("R8").Value2

# 2) Overwrite row 8 with the new weekly figures.
$ws.Range("D8").Value = 44448
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 580
